{"js": "// Resume edits: update Sling TV company description, rewrite/trim the Vail\n// Resorts accomplishment bullets (dropping two of them), and swap two\n// entries in the Skills word-jumble table.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Helper: find the (first) paragraph whose text starts with the given\n// needle, so the script is resilient to small index shifts.\nfunction findParagraphIndex(items, needle) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(needle) === 0) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + needle);\n}\n\nconst items = paragraphs.items;\n\n// 1) Sling TV company description.\nconst idxSling = findParagraphIndex(\n  items,\n  \"Streaming service providing affordable live TV cable alternatives\"\n);\nitems[idxSling].insertText(\n  \"Subsidiary of DISH Network, a subscription-based streaming service providing affordable live TV cable alternatives\",\n  Word.InsertLocation.replace\n);\n\n// 2) Vail Resorts accomplishment bullets. Locate the block by its first\n// bullet, then operate on the fixed-size group that follows it.\nconst idxLead = findParagraphIndex(\n  items,\n  \"4+ years as lead developer responsible for translating marketing strategy\"\n);\n\nconst pLead = items[idxLead];           // -> rewritten (shortened)\nconst pCampaigns = items[idxLead + 1];  // \"Designed and automated...\"   -> rewritten\nconst pMerger = items[idxLead + 2];     // \"Unit, integration, system...\" -> rewritten\nconst pOptimized = items[idxLead + 3];  // \"Optimized 30M+...\"            -> removed\nconst pTools = items[idxLead + 4];      // \"Created and maintained...\"    -> rewritten (shortened)\nconst pFacebook = items[idxLead + 5];   // \"Managed Facebook ad...\"       -> removed\n\npLead.insertText(\n  \"4+ years as lead for translating marketing strategy into executable code for Season Pass email communications.\",\n  Word.InsertLocation.replace\n);\n\npCampaigns.insertText(\n  \"Unit, and acceptance testing for integration of 7+ separate databases from acquired companies using SQL.\",\n  Word.InsertLocation.replace\n);\n\npMerger.insertText(\n  \"Optimized 30M+ customer data records according to email marketing industry best practices to boost deliverability and overall revenue.\",\n  Word.InsertLocation.replace\n);\n\npTools.insertText(\n  \"Created and maintained custom productivity and reporting tools in Python and Alteryx.\",\n  Word.InsertLocation.replace\n);\n\n// Remove the two bullets that no longer exist in the revised resume.\npOptimized.delete();\npFacebook.delete();\n\nawait context.sync();\n\n// 3) Skills table: \"Data Analysis\" -> \"Tableau\", \"Data Visualization\" -> \"Data Analysis\".\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst skillsTable = tables.items[0];\n\nconst rows = skillsTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row index 2 (0-based) is Git / Data Analysis / Data Visualization / Marketing Automation.\nconst thirdRow = skillsTable.rows.items[2];\nconst rowCells = thirdRow.cells;\nrowCells.load(\"items\");\nawait context.sync();\n\nconst cellDataAnalysis = rowCells.items[1].body.paragraphs.getFirst();\ncellDataAnalysis.insertText(\"Tableau\", Word.InsertLocation.replace);\n\nconst cellDataVisualization = rowCells.items[2].body.paragraphs.getFirst();\ncellDataVisualization.insertText(\"Data Analysis\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// 4) Tighten the spacing-after on the ResumeWordJumble style (6pt -> 4pt,\n// i.e. 120 -> 80 twips), matching the resaved style definition.\nconst jumbleStyle = context.document.getStyles().getByNameOrNullObject(\"ResumeWordJumble\");\njumbleStyle.paragraphFormat.spaceAfter = 4;\nawait context.sync();\n", "ps1": "# Resume edits: update Sling TV company description, rewrite/trim the Vail\n# Resorts accomplishment bullets (dropping two of them), swap two entries in\n# the Skills word-jumble table, and tighten the ResumeWordJumble spacing.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) Sling TV company description.\nReplace-ExactText `\n    \"Streaming service providing affordable live TV cable alternatives\" `\n    \"Subsidiary of DISH Network, a subscription-based streaming service providing affordable live TV cable alternatives\"\n\n# 2) Vail Resorts accomplishment bullets. Remove the two bullets that no\n# longer exist in the revised resume FIRST (while their text is still\n# unique), so later text rewrites can't create duplicate-text ambiguity.\nfunction Remove-ParagraphStartingWith($needle) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.StartsWith($needle)) {\n            $p.Range.Delete()\n            return\n        }\n    }\n}\n\nRemove-ParagraphStartingWith \"Optimized 30M+ customer data records according to email marketing industry best practices to boost deliverability and overall revenue.\"\nRemove-ParagraphStartingWith \"Managed Facebook ad deployment, pulling segments via enterprise database to upload as Custom Audiences.\"\n\nReplace-ExactText `\n    \"4+ years as lead developer responsible for translating marketing strategy into executable code for Season Pass email communications, the foremost revenue driver for this Fortune 1000 company.\" `\n    \"4+ years as lead for translating marketing strategy into executable code for Season Pass email communications.\"\n\nReplace-ExactText `\n    \"Designed and automated hundreds of marketing campaigns, each with unique strategy and data requirements, often sending to millions of customers with scores of personalized content segments.\" `\n    \"Unit, and acceptance testing for integration of 7+ separate databases from acquired companies using SQL.\"\n\nReplace-ExactText `\n    \"Unit, integration, system, and acceptance testing the merger of 7+ separate databases from acquired companies using SQL and Alteryx.\" `\n    \"Optimized 30M+ customer data records according to email marketing industry best practices to boost deliverability and overall revenue.\"\n\nReplace-ExactText `\n    \"Created and maintained custom productivity and reporting tools in Python and Alteryx used by multiple team members and visible across the organization.\" `\n    \"Created and maintained custom productivity and reporting tools in Python and Alteryx.\"\n\n# 3) Skills table: \"Data Analysis\" -> \"Tableau\", \"Data Visualization\" -> \"Data Analysis\".\n$skillsTable = $d.Tables.Item(1)\n$skillsTable.Cell(3, 2).Range.Text = \"Tableau\"\n$skillsTable.Cell(3, 3).Range.Text = \"Data Analysis\"\n\n# 4) Tighten the spacing-after on the ResumeWordJumble style (6pt -> 4pt,\n# i.e. 120 -> 80 twips), matching the resaved style definition.\n$jumbleStyle = $d.Styles.Item(\"ResumeWordJumble\")\n$jumbleStyle.ParagraphFormat.SpaceAfter = 4\n"}
